$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2023-06-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-23 Friday", 2) | Out-Null

# Update table cell values by position (row, col)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "79-37=42"
$tbl.Cell(1, 2).Range.Text = "14+85=99"
$tbl.Cell(1, 3).Range.Text = "33+0=33"
$tbl.Cell(1, 4).Range.Text = "85-57=28"
$tbl.Cell(1, 5).Range.Text = "75-20=55"
$tbl.Cell(2, 1).Range.Text = "84-57=27"
$tbl.Cell(2, 2).Range.Text = "43-28=15"
$tbl.Cell(2, 3).Range.Text = "15+80=95"
$tbl.Cell(2, 4).Range.Text = "73-57=16"
$tbl.Cell(2, 5).Range.Text = "77-60=17"
$tbl.Cell(3, 1).Range.Text = "30-4=26"
$tbl.Cell(3, 2).Range.Text = "28+39=67"
$tbl.Cell(3, 3).Range.Text = "1-0=1"
$tbl.Cell(3, 4).Range.Text = "54-5=49"
$tbl.Cell(3, 5).Range.Text = "78-40=38"
$tbl.Cell(4, 1).Range.Text = "93-75=18"
$tbl.Cell(4, 2).Range.Text = "64-24=40"
$tbl.Cell(4, 3).Range.Text = "70-60=10"
$tbl.Cell(4, 4).Range.Text = "21-7=14"
$tbl.Cell(4, 5).Range.Text = "49+29=78"
$tbl.Cell(5, 1).Range.Text = "7+20=27"
$tbl.Cell(5, 2).Range.Text = "72-42=30"
$tbl.Cell(5, 3).Range.Text = "68+23=91"
$tbl.Cell(5, 4).Range.Text = "36+63=99"
$tbl.Cell(5, 5).Range.Text = "94-93=1"
$tbl.Cell(6, 1).Range.Text = "18+72=90"
$tbl.Cell(6, 2).Range.Text = "56+2=58"
$tbl.Cell(6, 3).Range.Text = "99-38=61"
$tbl.Cell(6, 4).Range.Text = "63-42=21"
$tbl.Cell(6, 5).Range.Text = "35+25=60"
$tbl.Cell(7, 1).Range.Text = "61+21=82"
$tbl.Cell(7, 2).Range.Text = "15+23=38"
$tbl.Cell(7, 3).Range.Text = "27+37=64"
$tbl.Cell(7, 4).Range.Text = "27-26=1"
$tbl.Cell(7, 5).Range.Text = "8+65=73"
$tbl.Cell(8, 1).Range.Text = "86-76=10"
$tbl.Cell(8, 2).Range.Text = "49-20=29"
$tbl.Cell(8, 3).Range.Text = "0+43=43"
$tbl.Cell(8, 4).Range.Text = "83-4=79"
$tbl.Cell(8, 5).Range.Text = "2+3=5"
$tbl.Cell(9, 1).Range.Text = "80-77=3"
$tbl.Cell(9, 2).Range.Text = "9+4=13"
$tbl.Cell(9, 3).Range.Text = "21-9=12"
$tbl.Cell(9, 4).Range.Text = "58+11=69"
$tbl.Cell(9, 5).Range.Text = "70+15=85"
$tbl.Cell(10, 1).Range.Text = "74+24=98"
$tbl.Cell(10, 2).Range.Text = "82-15=67"
$tbl.Cell(10, 3).Range.Text = "0+3=3"
$tbl.Cell(10, 4).Range.Text = "46-40=6"
$tbl.Cell(10, 5).Range.Text = "23+59=82"
$tbl.Cell(11, 1).Range.Text = "74-62=12"
$tbl.Cell(11, 2).Range.Text = "40+6=46"
$tbl.Cell(11, 3).Range.Text = "89+7=96"
$tbl.Cell(11, 4).Range.Text = "24+66=90"
$tbl.Cell(11, 5).Range.Text = "8+67=75"
$tbl.Cell(12, 1).Range.Text = "78-39=39"
$tbl.Cell(12, 2).Range.Text = "67-0=67"
$tbl.Cell(12, 3).Range.Text = "82-11=71"
$tbl.Cell(12, 4).Range.Text = "31+23=54"
$tbl.Cell(12, 5).Range.Text = "67+11=78"
$tbl.Cell(13, 1).Range.Text = "17+47=64"
$tbl.Cell(13, 2).Range.Text = "92-6=86"
$tbl.Cell(13, 3).Range.Text = "70-21=49"
$tbl.Cell(13, 4).Range.Text = "31+46=77"
$tbl.Cell(13, 5).Range.Text = "3-1=2"
$tbl.Cell(14, 1).Range.Text = "68+26=94"
$tbl.Cell(14, 2).Range.Text = "11+25=36"
$tbl.Cell(14, 3).Range.Text = "67-16=51"
$tbl.Cell(14, 4).Range.Text = "39-11=28"
$tbl.Cell(14, 5).Range.Text = "41+34=75"
$tbl.Cell(15, 1).Range.Text = "20+26=46"
$tbl.Cell(15, 2).Range.Text = "60-37=23"
$tbl.Cell(15, 3).Range.Text = "44+11=55"
$tbl.Cell(15, 4).Range.Text = "37+17=54"
$tbl.Cell(15, 5).Range.Text = "15+10=25"
$tbl.Cell(16, 1).Range.Text = "75-26=49"
$tbl.Cell(16, 2).Range.Text = "28+55=83"
$tbl.Cell(16, 3).Range.Text = "77-69=8"
$tbl.Cell(16, 4).Range.Text = "93-44=49"
$tbl.Cell(16, 5).Range.Text = "67-20=47"
$tbl.Cell(17, 1).Range.Text = "33-15=18"
$tbl.Cell(17, 2).Range.Text = "3+65=68"
$tbl.Cell(17, 3).Range.Text = "15+18=33"
$tbl.Cell(17, 4).Range.Text = "38+6=44"
$tbl.Cell(17, 5).Range.Text = "12-0=12"
$tbl.Cell(18, 1).Range.Text = "45+35=80"
$tbl.Cell(18, 2).Range.Text = "48-37=11"
$tbl.Cell(18, 3).Range.Text = "6+37=43"
$tbl.Cell(18, 4).Range.Text = "91-29=62"
$tbl.Cell(18, 5).Range.Text = "6+9=15"
$tbl.Cell(19, 1).Range.Text = "44+5=49"
$tbl.Cell(19, 2).Range.Text = "36-19=17"
$tbl.Cell(19, 3).Range.Text = "68+14=82"
$tbl.Cell(19, 4).Range.Text = "32+7=39"
$tbl.Cell(19, 5).Range.Text = "58+23=81"
$tbl.Cell(20, 1).Range.Text = "49-17=32"
$tbl.Cell(20, 2).Range.Text = "73-51=22"
$tbl.Cell(20, 3).Range.Text = "26-0=26"
$tbl.Cell(20, 4).Range.Text = "45+19=64"
$tbl.Cell(20, 5).Range.Text = "73-41=32"

Write-Output "done"